# Update CDA Logical model metadata for ST.r2b
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Version value (row 3)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update the Date value (row 8)
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row after "Contact" (row 10), before "Description" (row 11)
$ws.Rows.Item(10).Copy()
$ws.Rows.Item(11).Insert()

# Copy the formatting of the row above down onto the newly inserted row
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

$excel.CutCopyMode = 0
